{"js": "// Replace each paragraph's text in document order (including table-cell\n// paragraphs) per the recorded diff: the date line plus the 25 division\n// problems in the practice table all get new values, one-to-one in order.\nconst replacements = [\n  \"2024-08-29 Thursday\",\n  \"60\u00f78=7, 4\",\n  \"13\u00f78=1, 5\",\n  \"38\u00f75=7, 3\",\n  \"64\u00f74=16, 0\",\n  \"39\u00f76=6, 3\",\n  \"20\u00f77=2, 6\",\n  \"70\u00f78=8, 6\",\n  \"74\u00f77=10, 4\",\n  \"72\u00f75=14, 2\",\n  \"19\u00f75=3, 4\",\n  \"84\u00f78=10, 4\",\n  \"29\u00f75=5, 4\",\n  \"74\u00f78=9, 2\",\n  \"39\u00f76=6, 3\",\n  \"59\u00f73=19, 2\",\n  \"97\u00f72=48, 1\",\n  \"30\u00f76=5, 0\",\n  \"33\u00f74=8, 1\",\n  \"74\u00f72=37, 0\",\n  \"28\u00f72=14, 0\",\n  \"54\u00f77=7, 5\",\n  \"47\u00f75=9, 2\",\n  \"28\u00f72=14, 0\",\n  \"80\u00f76=13, 2\",\n  \"87\u00f76=14, 3\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text && text.trim().length > 0) {\n    para.insertText(replacements[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date and every division-problem answer cell.\n# All \"before\" strings are unique in the document, so a literal\n# Find/Replace (no wildcards) for each pair safely retargets only the\n# intended run while leaving its formatting (rFonts/sz) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-08-28 Wednesday\", \"2024-08-29 Thursday\"),\n    @(\"42\u00f72=21, 0\", \"60\u00f78=7, 4\"),\n    @(\"12\u00f73=4, 0\", \"13\u00f78=1, 5\"),\n    @(\"60\u00f72=30, 0\", \"38\u00f75=7, 3\"),\n    @(\"27\u00f77=3, 6\", \"64\u00f74=16, 0\"),\n    @(\"18\u00f78=2, 2\", \"39\u00f76=6, 3\"),\n    @(\"16\u00f72=8, 0\", \"20\u00f77=2, 6\"),\n    @(\"69\u00f74=17, 1\", \"70\u00f78=8, 6\"),\n    @(\"45\u00f76=7, 3\", \"74\u00f77=10, 4\"),\n    @(\"35\u00f73=11, 2\", \"72\u00f75=14, 2\"),\n    @(\"90\u00f72=45, 0\", \"19\u00f75=3, 4\"),\n    @(\"87\u00f77=12, 3\", \"84\u00f78=10, 4\"),\n    @(\"95\u00f73=31, 2\", \"29\u00f75=5, 4\"),\n    @(\"23\u00f75=4, 3\", \"74\u00f78=9, 2\"),\n    @(\"36\u00f74=9, 0\", \"39\u00f76=6, 3\"),\n    @(\"91\u00f74=22, 3\", \"59\u00f73=19, 2\"),\n    @(\"54\u00f78=6, 6\", \"97\u00f72=48, 1\"),\n    @(\"29\u00f72=14, 1\", \"30\u00f76=5, 0\"),\n    @(\"70\u00f76=11, 4\", \"33\u00f74=8, 1\"),\n    @(\"22\u00f78=2, 6\", \"74\u00f72=37, 0\"),\n    @(\"80\u00f77=11, 3\", \"28\u00f72=14, 0\"),\n    @(\"20\u00f72=10, 0\", \"54\u00f77=7, 5\"),\n    @(\"85\u00f75=17, 0\", \"47\u00f75=9, 2\"),\n    @(\"10\u00f76=1, 4\", \"28\u00f72=14, 0\"),\n    @(\"20\u00f74=5, 0\", \"80\u00f76=13, 2\"),\n    @(\"25\u00f75=5, 0\", \"87\u00f76=14, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
